$d = $word.ActiveDocument
$xml = $d.WordOpenXML

# 1) Add the xmlns:oel="http://schemas.microsoft.com/office/2019/extlst" namespace
#    declaration to the root elements of word/document.xml, word/endnotes.xml,
#    word/footnotes.xml and word/numbering.xml (inserted right after the
#    existing xmlns:o="urn:schemas-microsoft-com:office:office" declaration).
#    word/settings.xml also declares xmlns:o but must stay untouched, so the
#    pattern is scoped to the specific root element names.
$pattern = '(<w:(?:document|endnotes|footnotes|numbering)\b[^>]*?xmlns:o="urn:schemas-microsoft-com:office:office" )(xmlns:r=)'
$replacement = '${1}xmlns:oel="http://schemas.microsoft.com/office/2019/extlst" ${2}'
$xml = [regex]::Replace($xml, $pattern, $replacement)

# 2) Stamp each <w:num> entry in word/numbering.xml with its durable numbering id.
$xml = $xml.Replace('<w:num w:numId="1">', '<w:num w:numId="1" w16cid:durableId="1741558365">')
$xml = $xml.Replace('<w:num w:numId="2">', '<w:num w:numId="2" w16cid:durableId="1740244445">')
$xml = $xml.Replace('<w:num w:numId="3">', '<w:num w:numId="3" w16cid:durableId="723135621">')
$xml = $xml.Replace('<w:num w:numId="4">', '<w:num w:numId="4" w16cid:durableId="1163856901">')
$xml = $xml.Replace('<w:num w:numId="5">', '<w:num w:numId="5" w16cid:durableId="1794473006">')
$xml = $xml.Replace('<w:num w:numId="6">', '<w:num w:numId="6" w16cid:durableId="180823632">')
$xml = $xml.Replace('<w:num w:numId="7">', '<w:num w:numId="7" w16cid:durableId="1275941500">')
$xml = $xml.Replace('<w:num w:numId="8">', '<w:num w:numId="8" w16cid:durableId="1201430537">')
$xml = $xml.Replace('<w:num w:numId="9">', '<w:num w:numId="9" w16cid:durableId="488205697">')
$xml = $xml.Replace('<w:num w:numId="10">', '<w:num w:numId="10" w16cid:durableId="725177963">')
$xml = $xml.Replace('<w:num w:numId="11">', '<w:num w:numId="11" w16cid:durableId="1861699342">')
$xml = $xml.Replace('<w:num w:numId="12">', '<w:num w:numId="12" w16cid:durableId="722489282">')
$xml = $xml.Replace('<w:num w:numId="13">', '<w:num w:numId="13" w16cid:durableId="1178231659">')

# 3) word/styles.xml: the "TOC Heading" (TOCHeading) style gains an rsid and an
#    explicit "no numbering" override (<w:numPr><w:numId w:val="0"/></w:numPr>)
#    at the top of its <w:pPr>.
$oldStyle = '<w:uiPriority w:val="39"/><w:unhideWhenUsed/><w:qFormat/><w:pPr><w:spacing w:before="240" w:line="259" w:lineRule="auto"/><w:outlineLvl w:val="9"/></w:pPr>'
$newStyle = '<w:uiPriority w:val="39"/><w:unhideWhenUsed/><w:qFormat/><w:rsid w:val="003E7207"/><w:pPr><w:numPr><w:numId w:val="0"/></w:numPr><w:spacing w:before="240" w:line="259" w:lineRule="auto"/><w:outlineLvl w:val="9"/></w:pPr>'
$xml = $xml.Replace($oldStyle, $newStyle)

$d.WordOpenXML = $xml
